# Updated symbol list on Sun Feb  5 22:38:37 UTC 2023 with GitHub Actions
#
# Refreshes the crypto-price snapshot: columns D (Price) and E
# (Volume(1h) % change) for the affected coin rows. Values are written
# with a leading apostrophe so Excel stores them as literal text
# (matching the existing inlineStr/text cells) instead of re-typing
# them as numbers or percentages, then the cell style is reset to
# "Normal" so no stray NumberFormat/quote-prefix style sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.08%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'43.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'5.44%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.560"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.37%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08097"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-3.93%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.673"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.61%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.284"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-4.79%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.888"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.37%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D10").Value = "'0.9362"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.02%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1169"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-5.97%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1896"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.54%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09588"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.75%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.04152"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'4.94%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1069"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.64%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001270"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.99%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005909"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-3.32%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.563"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.71%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.69%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.523"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-6.90%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1364"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.02%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2588"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.01%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04327"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.03%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001239"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.59%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004389"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.43%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001228"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.05%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.04%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02659"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-6.07%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05488"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.99%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.01145"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'27.78%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007651"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.48%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1397"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.94%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002104"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.02%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009635"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-7.17%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007005"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.15%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.28%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003546"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'10.43%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002272"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.37%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.28%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.28%"
$ws.Range("E51").Style = "Normal"
